# Apply update to "Database_Thresholds_details" worksheet:
#  1. Update the "ScriptLatestRunVersion" text (Git Commit ID) for all data rows.
#  2. Update the "pid" value in column AH for all data rows from 30656 to 25596.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldCommitText = "IndicatorQuantiles.R, Git Commit ID: 54e4488a188edf59eafc4b9cfe53dc7125db7b32"
$newCommitText = "IndicatorQuantiles.R, Git Commit ID: 0e4152332be22faf035a2e2fc83ad2cca4c8a7fc"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 80
}

for ($row = 2; $row -le $lastRow; $row++) {
    $ajCell = $ws.Cells.Item($row, 36) # column AJ
    if ($ajCell.Value2 -eq $oldCommitText) {
        $ajCell.Value2 = $newCommitText
    }

    $ahCell = $ws.Cells.Item($row, 34) # column AH
    if ($ahCell.Value2 -eq 30656) {
        $ahCell.Value2 = 25596
    }
}
